{"js": "// Office.js (Word JavaScript API) edit script.\n// Updates the heading date and every math-problem cell in the single\n// table to the values from the target revision.\n\nconst NEW_DATE = \"2025-08-14 Thursday\";\nconst NEW_VALUES = [\n  [\"72-36=\", \"9+22=\", \"49+14=\", \"62-58=\", \"8+69=\"],\n  [\"9+86=\", \"57+34=\", \"7+8=\", \"34-17=\", \"35+36=\"],\n  [\"97-88=\", \"35+39=\", \"45+39=\", \"68+18=\", \"27+26=\"],\n  [\"96-28=\", \"19+5=\", \"95-57=\", \"27+38=\", \"72-37=\"],\n  [\"35-9=\", \"90-42=\", \"91-65=\", \"2+89=\", \"90-76=\"],\n  [\"92-55=\", \"14-6=\", \"17+75=\", \"39+6=\", \"15+9=\"],\n  [\"56-49=\", \"95-86=\", \"41-2=\", \"93-48=\", \"44+19=\"],\n  [\"91-75=\", \"34+59=\", \"28+8=\", \"48+8=\", \"25+58=\"],\n  [\"16+78=\", \"81-57=\", \"83-24=\", \"93-18=\", \"7+47=\"],\n  [\"89+5=\", \"44+7=\", \"95-39=\", \"55+27=\", \"17+78=\"],\n  [\"32-8=\", \"43+29=\", \"25+59=\", \"18+53=\", \"69+3=\"],\n  [\"86-9=\", \"73-65=\", \"9+48=\", \"18+58=\", \"87-79=\"],\n  [\"26+48=\", \"28+69=\", \"18+69=\", \"65-7=\", \"57-38=\"],\n  [\"49+2=\", \"83-17=\", \"27+25=\", \"86+9=\", \"2+19=\"],\n  [\"48+4=\", \"77+17=\", \"49+22=\", \"68+4=\", \"86-17=\"],\n  [\"15+67=\", \"54-16=\", \"27+29=\", \"84-78=\", \"45+28=\"],\n  [\"73-44=\", \"81-54=\", \"37+19=\", \"37+5=\", \"43-28=\"],\n  [\"68+14=\", \"68+24=\", \"43+8=\", \"9+25=\", \"92-69=\"],\n  [\"94-76=\", \"23-19=\", \"17+25=\", \"90-5=\", \"19+18=\"],\n  [\"26+18=\", \"70-21=\", \"93-58=\", \"5+26=\", \"51-28=\"]\n];\n\n// --- Update the date heading (first paragraph, before the table) ---\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.insertText(NEW_DATE, Word.InsertLocation.replace);\n\n// --- Update every cell of the (only) table in the document ---\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nfor (let r = 0; r < NEW_VALUES.length; r++) {\n  const row = NEW_VALUES[r];\n  for (let c = 0; c < row.length; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = row[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Updates the heading date and every math-problem cell in the single\n# table to the values from the target revision.\n\n$NewDate = \"2025-08-14 Thursday\"\n$NewValues = @(\n    @(\"72-36=\", \"9+22=\", \"49+14=\", \"62-58=\", \"8+69=\"),\n    @(\"9+86=\", \"57+34=\", \"7+8=\", \"34-17=\", \"35+36=\"),\n    @(\"97-88=\", \"35+39=\", \"45+39=\", \"68+18=\", \"27+26=\"),\n    @(\"96-28=\", \"19+5=\", \"95-57=\", \"27+38=\", \"72-37=\"),\n    @(\"35-9=\", \"90-42=\", \"91-65=\", \"2+89=\", \"90-76=\"),\n    @(\"92-55=\", \"14-6=\", \"17+75=\", \"39+6=\", \"15+9=\"),\n    @(\"56-49=\", \"95-86=\", \"41-2=\", \"93-48=\", \"44+19=\"),\n    @(\"91-75=\", \"34+59=\", \"28+8=\", \"48+8=\", \"25+58=\"),\n    @(\"16+78=\", \"81-57=\", \"83-24=\", \"93-18=\", \"7+47=\"),\n    @(\"89+5=\", \"44+7=\", \"95-39=\", \"55+27=\", \"17+78=\"),\n    @(\"32-8=\", \"43+29=\", \"25+59=\", \"18+53=\", \"69+3=\"),\n    @(\"86-9=\", \"73-65=\", \"9+48=\", \"18+58=\", \"87-79=\"),\n    @(\"26+48=\", \"28+69=\", \"18+69=\", \"65-7=\", \"57-38=\"),\n    @(\"49+2=\", \"83-17=\", \"27+25=\", \"86+9=\", \"2+19=\"),\n    @(\"48+4=\", \"77+17=\", \"49+22=\", \"68+4=\", \"86-17=\"),\n    @(\"15+67=\", \"54-16=\", \"27+29=\", \"84-78=\", \"45+28=\"),\n    @(\"73-44=\", \"81-54=\", \"37+19=\", \"37+5=\", \"43-28=\"),\n    @(\"68+14=\", \"68+24=\", \"43+8=\", \"9+25=\", \"92-69=\"),\n    @(\"94-76=\", \"23-19=\", \"17+25=\", \"90-5=\", \"19+18=\"),\n    @(\"26+18=\", \"70-21=\", \"93-58=\", \"5+26=\", \"51-28=\")\n)\n\n$d = $word.ActiveDocument\n\n# --- Update the date heading (first paragraph, before the table) ---\n$dateParagraph = $d.Paragraphs.Item(1)\n$dateParagraph.Range.Text = $NewDate\n\n# --- Update every cell of the (only) table in the document ---\n$t = $d.Tables.Item(1)\nfor ($r = 1; $r -le $NewValues.Count; $r++) {\n    $row = $NewValues[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $row[$c - 1]\n    }\n}\n\n"}
